# 自动更新Excel文件
# Recomputes column E ("剩余" / days remaining) for each shop row from
# column D ("总天" / contract length in days) and column F ("开始时间" /
# start date, stored as an 8-digit yyyymmdd number), as of "today".
#
# Each row's end date = F + D days. Remaining = end date - today.
# If a contract's remaining days would reach zero or below, it is treated
# as renewed starting today: F is reset to today and E is reset to the
# full contract length D.
#
# "Today" advanced by one day since the sheet was last refreshed, so every
# remaining-days count drops by 1 (unless the contract rolled over, which
# resets it to D with F = today).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Julian Day Number helpers (pure integer/date-string arithmetic; no
# ---- reliance on System.DateTime members, which aren't reliably bridged
# ---- in this COM host). ----

function YmdStrToJdn($ymdText) {
    $yPart = [int]$ymdText.Substring(0, 4)
    $mPart = [int]$ymdText.Substring(4, 2)
    $dPart = [int]$ymdText.Substring(6, 2)
    $adj = [math]::Floor((14 - $mPart) / 12)
    $yAdj = $yPart + 4800 - $adj
    $mAdj = $mPart + 12 * $adj - 3
    $jdnOut = $dPart + [math]::Floor((153 * $mAdj + 2) / 5) + 365 * $yAdj + [math]::Floor($yAdj / 4) - [math]::Floor($yAdj / 100) + [math]::Floor($yAdj / 400) - 32045
    return $jdnOut
}

function JdnToYmdStr($jdnIn) {
    $a2 = $jdnIn + 32044
    $b2 = [math]::Floor((4 * $a2 + 3) / 146097)
    $c2 = $a2 - [math]::Floor((146097 * $b2) / 4)
    $d2 = [math]::Floor((4 * $c2 + 3) / 1461)
    $e2 = $c2 - [math]::Floor((1461 * $d2) / 4)
    $m2 = [math]::Floor((5 * $e2 + 2) / 153)
    $dayOut = $e2 - [math]::Floor((153 * $m2 + 2) / 5) + 1
    $monthOut = $m2 + 3 - 12 * [math]::Floor($m2 / 10)
    $yearOut = 100 * $b2 + $d2 - 4800 + [math]::Floor($m2 / 10)
    $ymdOut = "{0:D4}{1:D2}{2:D2}" -f $yearOut, $monthOut, $dayOut
    return $ymdOut
}

# The workbook's reference date moved forward one day (per the nightly
# refresh that produced this edit).
$todayStr = "20260219"
$todayJdn = YmdStrToJdn $todayStr

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $totalDays = $ws.Cells.Item($row, 4).Value2
    if ($totalDays -eq $null) {
        continue
    }

    $startRaw = $ws.Cells.Item($row, 6).Value2
    if ($startRaw -eq $null) {
        continue
    }
    $startStr = [string]$startRaw

    # Skip rows whose start date isn't a clean 8-digit yyyymmdd (data-entry
    # typos, e.g. row 36's "202510929") rather than mis-recompute them.
    if ($startStr.Length -ne 8) {
        continue
    }

    $startJdn = YmdStrToJdn $startStr
    $endJdn = $startJdn + $totalDays
    $remaining = $endJdn - $todayJdn

    if ($remaining -le 0) {
        # Contract lapsed - renew it starting today.
        $newRemaining = $totalDays
        $newStartStr = $todayStr
    } else {
        $newRemaining = $remaining
        $newStartStr = $startStr
    }

    $ws.Cells.Item($row, 5).Value = $newRemaining
    if ($newStartStr -ne $startStr) {
        $ws.Cells.Item($row, 6).Value = [double]$newStartStr
    }
}
